$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room for the new "DC Bias [V]" column -----------------------
# Before: A Board ID | B Assembled? | C C34/C35 | D R22/R24 | E R23/R25 | F Location | G Owner
# After:  A Board ID | B Assembled? | C C34/C35 | D R22/R24(...) | E R23/R25(...) | F DC Bias [V] | G Location | H Owner
$ws.Columns("F").Insert()

# --- Header row ---------------------------------------------------------
# Set E1 before D1 so the new shared strings land in the same table order
# as the authored workbook (R23/R25 text registered before R22/R24 text).
$ws.Range("E1").Value = "R23/R25 (Pullup to 1.8V)"
$ws.Range("D1").Value = "R22/R24 (Pulldown to GND)"
$ws.Range("F1").Value = "DC Bias [V]"
$ws.Range("F1").NumberFormat = "0.00"

# --- Row 2 (CLK1): add the measured pulldown/pullup resistor jumper photo data
$ws.Range("D2").Value = 615
$ws.Range("E2").Value = 690
$ws.Range("F2").Formula = "=1.8*D2/(D2+E2)"
$ws.Range("F2").NumberFormat = "0.00"

# --- Row 3 (CLK2): DNP (do-not-populate) jumper configuration
$ws.Range("C3").Value = "1nF"
$ws.Range("D3").Value = "DNP"
$ws.Range("E3").Value = "DNP"

# --- Best-effort column widths to mirror the widened header columns -----
$ws.Columns("D").ColumnWidth = 23.76
$ws.Columns("E").ColumnWidth = 20.92
$ws.Columns("F").ColumnWidth = 11.42

# --- Match the author's final cursor position ---------------------------
$null = $ws.Range("F4").Select()
